# Update "想去人数" (F column) counts on the 展览 and 全部类型 sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $ws.Range("F2").Value = 72
    $ws.Range("F3").Value = 1437
    $ws.Range("F5").Value = 19
    $ws.Range("F7").Value = 38
    $ws.Range("F8").Value = 218
}
